$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix the lat/lon "," <-> "." SUBSTITUTE() formulas in columns E and F.
#    The sheet previously turned "." into "," (European decimal style); it
#    now turns "," into "." instead. These are shared-formula ranges, so
#    setting the formula on the whole range keeps them shared with the
#    correct relative references.
# ---------------------------------------------------------------------------
$ws.Range("E2").Formula  = '=SUBSTITUTE( SUBSTITUTE(C2,"N",""),",",".")'
$ws.Range("F2").Formula  = '=SUBSTITUTE( SUBSTITUTE(D2,"E",""),",",".")'

$ws.Range("E3:E66").Formula   = '=SUBSTITUTE( SUBSTITUTE(C3,"N",""),",",".")'
$ws.Range("F3:F66").Formula   = '=SUBSTITUTE( SUBSTITUTE(D3,"E",""),",",".")'

$ws.Range("E67:E130").Formula = '=SUBSTITUTE( SUBSTITUTE(C67,"N",""),",",".")'
$ws.Range("F67:F130").Formula = '=SUBSTITUTE( SUBSTITUTE(D67,"E",""),",",".")'

$ws.Range("E131:E177").Formula = '=SUBSTITUTE( SUBSTITUTE(C131,"N",""),",",".")'
$ws.Range("F131:F177").Formula = '=SUBSTITUTE( SUBSTITUTE(D131,"E",""),",",".")'

# ---------------------------------------------------------------------------
# 2) Rename the MIL1/MIL2/MIL3/MIL5/MIL7/MIL9 stand labels to the
#    zero-padded MIL01/MIL02/MIL03/MIL05/MIL07/MIL09 form (MIL11/13/14/15
#    keep their existing labels/order).
# ---------------------------------------------------------------------------
$ws.Cells.Item(178, 2).Value = "MIL01"
$ws.Cells.Item(179, 2).Value = "MIL02"
$ws.Cells.Item(180, 2).Value = "MIL03"
$ws.Cells.Item(181, 2).Value = "MIL05"
$ws.Cells.Item(182, 2).Value = "MIL07"
$ws.Cells.Item(183, 2).Value = "MIL09"

# ---------------------------------------------------------------------------
# 3) Append the new "apron-GA" parking stands (GA01..GA10) as rows 188-197.
#    Values are entered in the same order the original author used (GA10
#    first, then GA01..GA09, then the "apron-GA" label column last) so the
#    shared-string table comes out in the same order as the source file.
# ---------------------------------------------------------------------------
$ws.Cells.Item(197, 2).Value = "GA10"
$ws.Cells.Item(188, 2).Value = "GA01"
$ws.Cells.Item(189, 2).Value = "GA02"
$ws.Cells.Item(190, 2).Value = "GA03"
$ws.Cells.Item(191, 2).Value = "GA04"
$ws.Cells.Item(192, 2).Value = "GA05"
$ws.Cells.Item(193, 2).Value = "GA06"
$ws.Cells.Item(194, 2).Value = "GA07"
$ws.Cells.Item(195, 2).Value = "GA08"
$ws.Cells.Item(196, 2).Value = "GA09"

for ($r = 188; $r -le 197; $r++) {
    $ws.Cells.Item($r, 1).Value = "apron-GA"
}

$jValues = @(50.896245766464403, 50.8963837586831,   50.895469514396801, 50.895766480337898, `
             50.895932323109399, 50.897002698050599, 50.8971752854549,   50.897315816901099, `
             50.8974619486347,   50.897637972154797)
$nValues = @(4.4672248103730601, 4.46778502646987,   4.4665204829185798, 4.4679075168651199, `
             4.4684262655108098, 4.4646363213440097, 4.4652878767248501, 4.46589393550867,   `
             4.4664137570469196, 4.4670373493481401)

$r = 188
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item($r, 10).Value = $jValues[$i]
    $ws.Cells.Item($r, 14).Value = $nValues[$i]
    $r++
}

# ---------------------------------------------------------------------------
# 4) Update the sheet selection to match the author's final view.
# ---------------------------------------------------------------------------
$ws.Range("C194").Select()
